$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.766.96'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.60'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7309'
$ws.Range("E5").Value = '  -1.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.75'
$ws.Range("E6").Value = '  -0.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9992'
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3128'
$ws.Range("E8").Value = '  -0.65%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07104'
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.32'
$ws.Range("E10").Value = '  -1.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08196'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7389'
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.328'
$ws.Range("E13").Value = '  -1.98%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.858.98'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.25'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.753.37'
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.010'
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '248.03'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007789'
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9975'
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.103.96'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.738'
$ws.Range("E24").Value = '  -2.95%  '
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.179'
$ws.Range("E26").Value = '  -1.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.59'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.50'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.006'
$ws.Range("E29").Value = '  -1.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.446'
$ws.Range("E30").Value = '  -1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.519'
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.519'
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.169'
$ws.Range("E33").Value = '  -2.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05289'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.230'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7423'
$ws.Range("E36").Value = '  -1.61%  '
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.691'
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01931'
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.733'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4446'
$ws.Range("E41").Value = '  -0.72%  '
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8678'
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.14'
$ws.Range("E44").Value = '  -1.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.044.12'
$ws.Range("E45").Value = '  -6.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9992'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.76'
$ws.Range("E47").Value = '  +0.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.435'
$ws.Range("E48").Value = '  -3.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.811'
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.539'
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.005.69'
$ws.Range("E51").Value = '  -0.19%  '
